$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 715
$ws.Cells.Item(715, 1).Value = 715
$ws.Cells.Item(715, 2).NumberFormat = "@"
$ws.Cells.Item(715, 2).Value = 'FCC'
$ws.Cells.Item(715, 3).NumberFormat = "@"
$ws.Cells.Item(715, 3).Value = 'TRT 3ª Região'
$ws.Cells.Item(715, 4).NumberFormat = "@"
$ws.Cells.Item(715, 4).Value = '2015'
$ws.Cells.Item(715, 5).NumberFormat = "@"
$ws.Cells.Item(715, 5).Value = 'A Responsabilidade Social está normatizada na ISO 26000 que versa pela incorporação de práticas socioambientais nos processos decisórios e a responsabilização pelos impactos de suas atividades na sociedade e meio ambiente. Para certificação a empresa deve integrar, implementar e promover'
$ws.Cells.Item(715, 6).NumberFormat = "@"
$ws.Cells.Item(715, 6).Value = 'Conhecimentos Específicos'
$ws.Cells.Item(715, 7).NumberFormat = "@"
$ws.Cells.Item(715, 7).Value = 'Sustentabilidade'
$ws.Cells.Item(715, 8).NumberFormat = "@"
$ws.Cells.Item(715, 8).Value = 'Médio'
$ws.Cells.Item(715, 9).NumberFormat = "@"
$ws.Cells.Item(715, 9).Value = 'ME'
$ws.Cells.Item(715, 10).NumberFormat = "@"
$ws.Cells.Item(715, 10).Value = 'as práticas econômicas a um comportamento socialmente responsável.'
$ws.Cells.Item(715, 11).NumberFormat = "@"
$ws.Cells.Item(715, 11).Value = 'o comportamento socialmente responsável em toda empresa.'
$ws.Cells.Item(715, 12).NumberFormat = "@"
$ws.Cells.Item(715, 12).Value = 'práticas socioeducativas nas atividades fim, visando a mitigação do impacto ambiental.'
$ws.Cells.Item(715, 13).NumberFormat = "@"
$ws.Cells.Item(715, 13).Value = 'um comportamento responsável no âmbito interno, visando uma melhor relação com o entorno.'
$ws.Cells.Item(715, 14).NumberFormat = "@"
$ws.Cells.Item(715, 14).Value = 'ações socioambientais sem impactar nos aspectos econômicos da empresa.'
$ws.Cells.Item(715, 15).NumberFormat = "@"
$ws.Cells.Item(715, 15).Value = 'B'
$ws.Cells.Item(715, 16).Value = 0
$ws.Cells.Item(715, 17).Value = 0

# Row 716
$ws.Cells.Item(716, 1).Value = 716
$ws.Cells.Item(716, 2).NumberFormat = "@"
$ws.Cells.Item(716, 2).Value = 'CESPE'
$ws.Cells.Item(716, 3).NumberFormat = "@"
$ws.Cells.Item(716, 3).Value = 'INSS'
$ws.Cells.Item(716, 4).NumberFormat = "@"
$ws.Cells.Item(716, 4).Value = '2008'
$ws.Cells.Item(716, 5).NumberFormat = "@"
$ws.Cells.Item(716, 5).Value = 'Em relação à responsabilidade socioambiental corporativa, julgue os itens seguintes.
A noção de responsabilidade socioambiental relaciona-se ao propósito de maximização dos lucros da corporação, que termina por beneficiar a sociedade na forma de empregos, salários e impostos, parte dos quais pode ser empregada em conservação ambiental.'
$ws.Cells.Item(716, 6).NumberFormat = "@"
$ws.Cells.Item(716, 6).Value = 'Conhecimentos Específicos'
$ws.Cells.Item(716, 7).NumberFormat = "@"
$ws.Cells.Item(716, 7).Value = 'Sustentabilidade'
$ws.Cells.Item(716, 8).NumberFormat = "@"
$ws.Cells.Item(716, 8).Value = 'Médio'
$ws.Cells.Item(716, 9).NumberFormat = "@"
$ws.Cells.Item(716, 9).Value = 'CE'
$ws.Cells.Item(716, 10).NumberFormat = "@"
$ws.Cells.Item(716, 10).Value = ''
$ws.Cells.Item(716, 11).NumberFormat = "@"
$ws.Cells.Item(716, 11).Value = ''
$ws.Cells.Item(716, 12).NumberFormat = "@"
$ws.Cells.Item(716, 12).Value = ''
$ws.Cells.Item(716, 13).NumberFormat = "@"
$ws.Cells.Item(716, 13).Value = ''
$ws.Cells.Item(716, 14).NumberFormat = "@"
$ws.Cells.Item(716, 14).Value = ''
$ws.Cells.Item(716, 15).NumberFormat = "@"
$ws.Cells.Item(716, 15).Value = 'E'
$ws.Cells.Item(716, 16).Value = 0
$ws.Cells.Item(716, 17).Value = 0

# Row 717
$ws.Cells.Item(717, 1).Value = 717
$ws.Cells.Item(717, 2).NumberFormat = "@"
$ws.Cells.Item(717, 2).Value = 'CESGRANRIO'
$ws.Cells.Item(717, 3).NumberFormat = "@"
$ws.Cells.Item(717, 3).Value = 'Petrobras'
$ws.Cells.Item(717, 4).NumberFormat = "@"
$ws.Cells.Item(717, 4).Value = '2012'
$ws.Cells.Item(717, 5).NumberFormat = "@"
$ws.Cells.Item(717, 5).Value = 'A NBR ISO 14001:2004 especifica os principais requisitos de um Sistema de Gestão Ambiental (SGA), de modo que as questões ambientais sejam integradas à administração global de uma organização ou empresa.
Qual é a denominação de um requisito do SGA e sua(s) respectiva(s) fase(s)?'
$ws.Cells.Item(717, 6).NumberFormat = "@"
$ws.Cells.Item(717, 6).Value = 'Conhecimentos Específicos'
$ws.Cells.Item(717, 7).NumberFormat = "@"
$ws.Cells.Item(717, 7).Value = 'Sustentabilidade'
$ws.Cells.Item(717, 8).NumberFormat = "@"
$ws.Cells.Item(717, 8).Value = 'Médio'
$ws.Cells.Item(717, 9).NumberFormat = "@"
$ws.Cells.Item(717, 9).Value = 'ME'
$ws.Cells.Item(717, 10).NumberFormat = "@"
$ws.Cells.Item(717, 10).Value = 'Auditoria do SGA e fase de análise crítica'
$ws.Cells.Item(717, 11).NumberFormat = "@"
$ws.Cells.Item(717, 11).Value = 'Documentação do SGA e fases de verificação e ação corretiva'
$ws.Cells.Item(717, 12).NumberFormat = "@"
$ws.Cells.Item(717, 12).Value = 'Monitoramento e medição e fase de planejamento'
$ws.Cells.Item(717, 13).NumberFormat = "@"
$ws.Cells.Item(717, 13).Value = 'Estrutura e responsabilidades e fase de política ambiental'
$ws.Cells.Item(717, 14).NumberFormat = "@"
$ws.Cells.Item(717, 14).Value = 'Preparação e atendimento a emergências e fases de implementação e operação'
$ws.Cells.Item(717, 15).NumberFormat = "@"
$ws.Cells.Item(717, 15).Value = 'E'
$ws.Cells.Item(717, 16).Value = 0
$ws.Cells.Item(717, 17).Value = 0

# Row 718
$ws.Cells.Item(718, 1).Value = 718
$ws.Cells.Item(718, 2).NumberFormat = "@"
$ws.Cells.Item(718, 2).Value = 'IFPA'
$ws.Cells.Item(718, 3).NumberFormat = "@"
$ws.Cells.Item(718, 3).Value = 'IFPA'
$ws.Cells.Item(718, 4).NumberFormat = "@"
$ws.Cells.Item(718, 4).Value = '2015'
$ws.Cells.Item(718, 5).NumberFormat = "@"
$ws.Cells.Item(718, 5).Value = 'Um sistema de gestão ambiental (SGA) se constitui em um conjunto de procedimentos sistematizados que são desenvolvidos para que as questões ambientais sejam integradas à administração global de um empreendimento e possibilite a obtenção de melhores resultados no desempenho global da empresa. Sendo assim, com relação aos elementos de um SGA, é incorreto afirmar que:'
$ws.Cells.Item(718, 6).NumberFormat = "@"
$ws.Cells.Item(718, 6).Value = 'Conhecimentos Específicos'
$ws.Cells.Item(718, 7).NumberFormat = "@"
$ws.Cells.Item(718, 7).Value = 'Sustentabilidade'
$ws.Cells.Item(718, 8).NumberFormat = "@"
$ws.Cells.Item(718, 8).Value = 'Médio'
$ws.Cells.Item(718, 9).NumberFormat = "@"
$ws.Cells.Item(718, 9).Value = 'ME'
$ws.Cells.Item(718, 10).NumberFormat = "@"
$ws.Cells.Item(718, 10).Value = 'Para que o comprometimento com a melhoria possa ser efetivo, todos os atores que constituem a organização devem de forma contínua aplicar a revisão do sistema de gestão ambiental, assegurando que este continue adequado e efetivo.'
$ws.Cells.Item(718, 11).NumberFormat = "@"
$ws.Cells.Item(718, 11).Value = 'A política ambiental dá um senso global de direção, apresenta os princípios de ação para uma organização, sendo estabelecidas metas relativas de desempenho e responsabilidade ambiental, contra as quais todas as ações subsequentes serão julgadas.'
$ws.Cells.Item(718, 12).NumberFormat = "@"
$ws.Cells.Item(718, 12).Value = 'Com base na política ambiental, a organização deve fazer um planejamento com o objetivo de atender aos requisitos estabelecidos.'
$ws.Cells.Item(718, 13).NumberFormat = "@"
$ws.Cells.Item(718, 13).Value = 'O processo de implementação e operação do SGA deve ser conduzido de forma a serem atingidos os objetivos e as metas estabelecidas.'
$ws.Cells.Item(718, 14).NumberFormat = "@"
$ws.Cells.Item(718, 14).Value = 'É necessário que sejam desenvolvidos procedimentos para monitorar e medir as principais características das operações e atividades que podem causar um impacto significativo no meio ambiente, ao mesmo tempo em que devem ser estabelecidos os procedimentos referentes às ações corretivas que devem ser tomadas para eliminar as causas reais ou potenciais, que poderiam resultar em um impacto no meio ambiente.'
$ws.Cells.Item(718, 15).NumberFormat = "@"
$ws.Cells.Item(718, 15).Value = 'A'
$ws.Cells.Item(718, 16).Value = 0
$ws.Cells.Item(718, 17).Value = 0

# Row 719
$ws.Cells.Item(719, 1).Value = 719
$ws.Cells.Item(719, 2).NumberFormat = "@"
$ws.Cells.Item(719, 2).Value = 'VUNESP'
$ws.Cells.Item(719, 3).NumberFormat = "@"
$ws.Cells.Item(719, 3).Value = 'CEAGESP'
$ws.Cells.Item(719, 4).NumberFormat = "@"
$ws.Cells.Item(719, 4).Value = '2010'
$ws.Cells.Item(719, 5).NumberFormat = "@"
$ws.Cells.Item(719, 5).Value = 'A necessidade de demonstrar, junto às partes interessadas, um comportamento ambiental aceitável, estimula as organizações a adotarem sistemas de gestão ambiental, como aquele prescrito na NBR ISO n.º 14001. De acordo com essa norma,'
$ws.Cells.Item(719, 6).NumberFormat = "@"
$ws.Cells.Item(719, 6).Value = 'Conhecimentos Específicos'
$ws.Cells.Item(719, 7).NumberFormat = "@"
$ws.Cells.Item(719, 7).Value = 'Sustentabilidade'
$ws.Cells.Item(719, 8).NumberFormat = "@"
$ws.Cells.Item(719, 8).Value = 'Médio'
$ws.Cells.Item(719, 9).NumberFormat = "@"
$ws.Cells.Item(719, 9).Value = 'ME'
$ws.Cells.Item(719, 10).NumberFormat = "@"
$ws.Cells.Item(719, 10).Value = 'as iniciativas da organização, em relação ao treinamento, conscientização e desenvolvimento de competências, além dos aspectos ambientais, devem contemplar os aspectos da segurança e saúde no trabalho, que impactam o sistema de produção.'
$ws.Cells.Item(719, 11).NumberFormat = "@"
$ws.Cells.Item(719, 11).Value = 'uma avaliação ambiental inicial deve cobrir quatro áreas: requisitos legais e regulamentares; identificação dos aspectos ambientais significativos; exame das práticas e procedimentos de gestão ambiental existentes e avaliação dos dados provenientes de investigações anteriores.'
$ws.Cells.Item(719, 12).NumberFormat = "@"
$ws.Cells.Item(719, 12).Value = 'para atender a seus objetivos com todas as partes interessadas em seu desempenho ambiental, a organização deve atender a seus requisitos legais, não sendo necessário subscrever eventuais códigos de prática da indústria, acordos voluntários e outras diretrizes de caráter não regulamentar.'
$ws.Cells.Item(719, 13).NumberFormat = "@"
$ws.Cells.Item(719, 13).Value = 'as auditorias não podem ser executadas por pessoal interno à organização e devem verificar de forma sistemática os documentos e as práticas para obter evidências de que a organização está em conformidade, de acordo com padrões objetivos, com o disposto na norma.'
$ws.Cells.Item(719, 14).NumberFormat = "@"
$ws.Cells.Item(719, 14).Value = 'todos os impactos ambientais possíveis nos recursos naturais, flora, fauna e suas interrelações devem receber adequada valoração de acordo com as diversas abordagens, como intensidade de reação da opinião pública e custo direto de medidas de mitigação.'
$ws.Cells.Item(719, 15).NumberFormat = "@"
$ws.Cells.Item(719, 15).Value = 'A'
$ws.Cells.Item(719, 16).Value = 0
$ws.Cells.Item(719, 17).Value = 0

# Row 720
$ws.Cells.Item(720, 1).Value = 720
$ws.Cells.Item(720, 2).NumberFormat = "@"
$ws.Cells.Item(720, 2).Value = 'Cespe'
$ws.Cells.Item(720, 3).NumberFormat = "@"
$ws.Cells.Item(720, 3).Value = 'Polícia Científica'
$ws.Cells.Item(720, 4).NumberFormat = "@"
$ws.Cells.Item(720, 4).Value = '2016'
$ws.Cells.Item(720, 5).NumberFormat = "@"
$ws.Cells.Item(720, 5).Value = 'Acerca da NBR ISO n.º 14.001, que estabelece diretrizes básicas para um sistema de gestão ambiental (SGA), assinale a opção correta.'
$ws.Cells.Item(720, 6).NumberFormat = "@"
$ws.Cells.Item(720, 6).Value = 'Conhecimentos Específicos'
$ws.Cells.Item(720, 7).NumberFormat = "@"
$ws.Cells.Item(720, 7).Value = 'Sustentabilidade'
$ws.Cells.Item(720, 8).NumberFormat = "@"
$ws.Cells.Item(720, 8).Value = 'Médio'
$ws.Cells.Item(720, 9).NumberFormat = "@"
$ws.Cells.Item(720, 9).Value = 'ME'
$ws.Cells.Item(720, 10).NumberFormat = "@"
$ws.Cells.Item(720, 10).Value = 'O SGA substitui as políticas de gestão ambiental anteriormente implementadas em uma organização.'
$ws.Cells.Item(720, 11).NumberFormat = "@"
$ws.Cells.Item(720, 11).Value = 'A política ambiental refere-se às intenções e aos princípios gerais da organização relacionados a seu desempenho ambiental, conforme formalmente expresso pela alta administração.'
$ws.Cells.Item(720, 12).NumberFormat = "@"
$ws.Cells.Item(720, 12).Value = 'Para que a implantação de um SGA em uma organização seja bem-sucedida, é necessário que as funções e responsabilidades ambientais sejam compreendidas como funções da gestão ambiental, ou seja, como atribuições dos seus gestores, de modo a não interferirem em outras áreas da organização, a fim de que conflitos internos sejam evitados.'
$ws.Cells.Item(720, 13).NumberFormat = "@"
$ws.Cells.Item(720, 13).Value = 'Um SGA implementado em determinada organização deve ser aplicado simultaneamente em todas as áreas de atividade dessa empresa, com o propósito de aprimorar o desempenho ambiental geral da organização.'
$ws.Cells.Item(720, 14).NumberFormat = "@"
$ws.Cells.Item(720, 14).Value = 'O aspecto ambiental de uma organização é um componente de suas atividades, seus produtos ou seus serviços que pode interagir com o meio ambiente sem causar impacto ambiental.'
$ws.Cells.Item(720, 15).NumberFormat = "@"
$ws.Cells.Item(720, 15).Value = 'B'
$ws.Cells.Item(720, 16).Value = 0
$ws.Cells.Item(720, 17).Value = 0

# Row 721
$ws.Cells.Item(721, 1).Value = 721
$ws.Cells.Item(721, 2).NumberFormat = "@"
$ws.Cells.Item(721, 2).Value = 'FCC'
$ws.Cells.Item(721, 3).NumberFormat = "@"
$ws.Cells.Item(721, 3).Value = 'PREFEITURA DE TERESINA - PI'
$ws.Cells.Item(721, 4).NumberFormat = "@"
$ws.Cells.Item(721, 4).Value = '2016'
$ws.Cells.Item(721, 5).NumberFormat = "@"
$ws.Cells.Item(721, 5).Value = 'A norma ISO 14001 define política ambiental como'
$ws.Cells.Item(721, 6).NumberFormat = "@"
$ws.Cells.Item(721, 6).Value = 'Conhecimentos Específicos'
$ws.Cells.Item(721, 7).NumberFormat = "@"
$ws.Cells.Item(721, 7).Value = 'Sustentabilidade'
$ws.Cells.Item(721, 8).NumberFormat = "@"
$ws.Cells.Item(721, 8).Value = 'Médio'
$ws.Cells.Item(721, 9).NumberFormat = "@"
$ws.Cells.Item(721, 9).Value = 'ME'
$ws.Cells.Item(721, 10).NumberFormat = "@"
$ws.Cells.Item(721, 10).Value = 'qualquer modificação do meio ambiente, adversa ou benéfica, que resulte, no todo ou em parte, dos aspectos ambientais da organização.'
$ws.Cells.Item(721, 11).NumberFormat = "@"
$ws.Cells.Item(721, 11).Value = 'requisito de desempenho detalhado, aplicável à organização ou à parte dela, resultante dos objetivos ambientais e que necessita ser estabelecido e atendido para que tais objetivos sejam atingidos.'
$ws.Cells.Item(721, 12).NumberFormat = "@"
$ws.Cells.Item(721, 12).Value = 'processo sistemático, independente e documentado para obter evidência e avaliá-la objetivamente para determinar a extensão na qual os critérios de auditoria do sistema da gestão ambiental estabelecidos pela organização são atendidos.'
$ws.Cells.Item(721, 13).NumberFormat = "@"
$ws.Cells.Item(721, 13).Value = 'elemento das atividades ou produtos ou serviços de uma organização que pode interagir com o meio ambiente e que pode causar impacto ambiental significativo.'
$ws.Cells.Item(721, 14).NumberFormat = "@"
$ws.Cells.Item(721, 14).Value = 'intenções e princípios gerais de uma organização em relação ao seu desempenho ambiental, conforme formalmente expresso pela alta administração.'
$ws.Cells.Item(721, 15).NumberFormat = "@"
$ws.Cells.Item(721, 15).Value = 'E'
$ws.Cells.Item(721, 16).Value = 0
$ws.Cells.Item(721, 17).Value = 0

# Row 722
$ws.Cells.Item(722, 1).Value = 722
$ws.Cells.Item(722, 2).NumberFormat = "@"
$ws.Cells.Item(722, 2).Value = 'FCC'
$ws.Cells.Item(722, 3).NumberFormat = "@"
$ws.Cells.Item(722, 3).Value = 'SEGEP-MA'
$ws.Cells.Item(722, 4).NumberFormat = "@"
$ws.Cells.Item(722, 4).Value = '2016'
$ws.Cells.Item(722, 5).NumberFormat = "@"
$ws.Cells.Item(722, 5).Value = 'O Sistema de Gestão Ambiental − SGA consiste na estrutura, responsabilidades, práticas, procedimentos, programas e recursos mobilizados para o atendimento da política ambiental de uma organização. É característica de um SGA:'
$ws.Cells.Item(722, 6).NumberFormat = "@"
$ws.Cells.Item(722, 6).Value = 'Conhecimentos Específicos'
$ws.Cells.Item(722, 7).NumberFormat = "@"
$ws.Cells.Item(722, 7).Value = 'Sustentabilidade'
$ws.Cells.Item(722, 8).NumberFormat = "@"
$ws.Cells.Item(722, 8).Value = 'Médio'
$ws.Cells.Item(722, 9).NumberFormat = "@"
$ws.Cells.Item(722, 9).Value = 'ME'
$ws.Cells.Item(722, 10).NumberFormat = "@"
$ws.Cells.Item(722, 10).Value = 'Garantir que o controle ambiental seja realizado pelo órgão ambiental oficial.'
$ws.Cells.Item(722, 11).NumberFormat = "@"
$ws.Cells.Item(722, 11).Value = 'Conter requisitos que podem ser auditados objetivamente para fins de certificação.'
$ws.Cells.Item(722, 12).NumberFormat = "@"
$ws.Cells.Item(722, 12).Value = 'A alta administração ser responsável pela implantação e manutenção de um sistema de gestão ambiental.'
$ws.Cells.Item(722, 13).NumberFormat = "@"
$ws.Cells.Item(722, 13).Value = 'Substituir o cumprimento dos requisitos legais e regulatórios.'
$ws.Cells.Item(722, 14).NumberFormat = "@"
$ws.Cells.Item(722, 14).Value = 'Não ser aplicável a todos os tipos e tamanhos de organizações.'
$ws.Cells.Item(722, 15).NumberFormat = "@"
$ws.Cells.Item(722, 15).Value = 'B'
$ws.Cells.Item(722, 16).Value = 0
$ws.Cells.Item(722, 17).Value = 0

# Row 723
$ws.Cells.Item(723, 1).Value = 723
$ws.Cells.Item(723, 2).NumberFormat = "@"
$ws.Cells.Item(723, 2).Value = 'IFRS'
$ws.Cells.Item(723, 3).NumberFormat = "@"
$ws.Cells.Item(723, 3).Value = 'IFRS'
$ws.Cells.Item(723, 4).NumberFormat = "@"
$ws.Cells.Item(723, 4).Value = '2016'
$ws.Cells.Item(723, 5).NumberFormat = "@"
$ws.Cells.Item(723, 5).Value = 'Em relação ao escopo da norma ABNT NBR ISO 14001:2015 são feitas as seguintes afirmativas.
I. Especifica os requisitos para um sistema de gestão ambiental que uma organização pode usar para aumentar seu desempenho ambiental.
II. Destina-se ao uso por uma organização que busca gerenciar suas responsabilidades ambientais de uma forma sistemática, que contribua para o pilar ambiental da sustentabilidade.
III. Auxilia a organização a alcançar os resultados pretendidos de seu sistema de gestão ambiental, os quais agreguem valor para o meio ambiente, à organização em si e suas partes interessadas.
IV. É aplicável a qualquer organização, independentemente do seu tamanho, tipo e natureza, e aplica-se aos aspectos ambientais das suas atividades, produtos e serviços que a organização determina poder controlar ou influenciar, considerando uma perspectiva de ciclo de vida.
V. Pode ser usada na íntegra ou em parte para sistematicamente melhorar a gestão ambiental. As declarações de conformidade, no entanto, só são aceitas se todos os requisitos forem incorporados ao sistema de gestão ambiental da organização e atendidos sem exclusões.
Assinale a alternativa em que todas as afirmativas estão CORRETAS:'
$ws.Cells.Item(723, 6).NumberFormat = "@"
$ws.Cells.Item(723, 6).Value = 'Conhecimentos Específicos'
$ws.Cells.Item(723, 7).NumberFormat = "@"
$ws.Cells.Item(723, 7).Value = 'Sustentabilidade'
$ws.Cells.Item(723, 8).NumberFormat = "@"
$ws.Cells.Item(723, 8).Value = 'Médio'
$ws.Cells.Item(723, 9).NumberFormat = "@"
$ws.Cells.Item(723, 9).Value = 'ME'
$ws.Cells.Item(723, 10).NumberFormat = "@"
$ws.Cells.Item(723, 10).Value = 'I, II, III, IV e V.'
$ws.Cells.Item(723, 11).NumberFormat = "@"
$ws.Cells.Item(723, 11).Value = 'Apenas I, II e V.'
$ws.Cells.Item(723, 12).NumberFormat = "@"
$ws.Cells.Item(723, 12).Value = 'Apenas I, II, III e IV.'
$ws.Cells.Item(723, 13).NumberFormat = "@"
$ws.Cells.Item(723, 13).Value = 'Apenas I, III, IV e V.'
$ws.Cells.Item(723, 14).NumberFormat = "@"
$ws.Cells.Item(723, 14).Value = 'Apenas II, III e IV.'
$ws.Cells.Item(723, 15).NumberFormat = "@"
$ws.Cells.Item(723, 15).Value = 'A'
$ws.Cells.Item(723, 16).Value = 0
$ws.Cells.Item(723, 17).Value = 0

# Row 724
$ws.Cells.Item(724, 1).Value = 724
$ws.Cells.Item(724, 2).NumberFormat = "@"
$ws.Cells.Item(724, 2).Value = 'Cesgranrio'
$ws.Cells.Item(724, 3).NumberFormat = "@"
$ws.Cells.Item(724, 3).Value = 'Petroquímica SUAPE'
$ws.Cells.Item(724, 4).NumberFormat = "@"
$ws.Cells.Item(724, 4).Value = '2011'
$ws.Cells.Item(724, 5).NumberFormat = "@"
$ws.Cells.Item(724, 5).Value = 'A ISO 14001 é a norma internacionalmente conhecida que apresenta um Sistema de Gestão Ambiental (SGA). O processo de implementação desse sistema possui quatro fases. A última fase desse processo é a de'
$ws.Cells.Item(724, 6).NumberFormat = "@"
$ws.Cells.Item(724, 6).Value = 'Conhecimentos Específicos'
$ws.Cells.Item(724, 7).NumberFormat = "@"
$ws.Cells.Item(724, 7).Value = 'Sustentabilidade'
$ws.Cells.Item(724, 8).NumberFormat = "@"
$ws.Cells.Item(724, 8).Value = 'Médio'
$ws.Cells.Item(724, 9).NumberFormat = "@"
$ws.Cells.Item(724, 9).Value = 'ME'
$ws.Cells.Item(724, 10).NumberFormat = "@"
$ws.Cells.Item(724, 10).Value = 'auditoria e certificação'
$ws.Cells.Item(724, 11).NumberFormat = "@"
$ws.Cells.Item(724, 11).Value = 'implantação das medidas de controle'
$ws.Cells.Item(724, 12).NumberFormat = "@"
$ws.Cells.Item(724, 12).Value = 'planejamento'
$ws.Cells.Item(724, 13).NumberFormat = "@"
$ws.Cells.Item(724, 13).Value = 'indicação dos desvios do processo'
$ws.Cells.Item(724, 14).NumberFormat = "@"
$ws.Cells.Item(724, 14).Value = 'monitoramento'
$ws.Cells.Item(724, 15).NumberFormat = "@"
$ws.Cells.Item(724, 15).Value = 'A'
$ws.Cells.Item(724, 16).Value = 0
$ws.Cells.Item(724, 17).Value = 0

# Row 725
$ws.Cells.Item(725, 1).Value = 725
$ws.Cells.Item(725, 2).NumberFormat = "@"
$ws.Cells.Item(725, 2).Value = 'UFSM'
$ws.Cells.Item(725, 3).NumberFormat = "@"
$ws.Cells.Item(725, 3).Value = 'UFSM'
$ws.Cells.Item(725, 4).NumberFormat = "@"
$ws.Cells.Item(725, 4).Value = '2018'
$ws.Cells.Item(725, 5).NumberFormat = "@"
$ws.Cells.Item(725, 5).Value = 'A base para um sistema de gestão ambiental é fundamentada no conceito Plan-Do-Check-Act (PDCA), um processo interativo utilizado pelas organizações para alcançara melhoria contínua.
Com relação ao ciclo PDCA, é INCORRETO afirmar:'
$ws.Cells.Item(725, 6).NumberFormat = "@"
$ws.Cells.Item(725, 6).Value = 'Conhecimentos Específicos'
$ws.Cells.Item(725, 7).NumberFormat = "@"
$ws.Cells.Item(725, 7).Value = 'Sustentabilidade'
$ws.Cells.Item(725, 8).NumberFormat = "@"
$ws.Cells.Item(725, 8).Value = 'Médio'
$ws.Cells.Item(725, 9).NumberFormat = "@"
$ws.Cells.Item(725, 9).Value = 'ME'
$ws.Cells.Item(725, 10).NumberFormat = "@"
$ws.Cells.Item(725, 10).Value = 'na etapa Plan, são definidos os objetivos ambientais que devem ser coerentes com a política ambiental, mensuráveis (se viável), monitorados, comunicados e atualizados, como apropriado.'
$ws.Cells.Item(725, 11).NumberFormat = "@"
$ws.Cells.Item(725, 11).Value = 'na etapa Act, poderão ser feitas melhorias como ação corretiva, melhoria contínua, mudança inovadora, inovação e reorganização.'
$ws.Cells.Item(725, 12).NumberFormat = "@"
$ws.Cells.Item(725, 12).Value = 'na etapa Check, os resultados são medidos em relação à política ambiental da organização, aos objetivos ambientais e a outro critério, usando indicadores.'
$ws.Cells.Item(725, 13).NumberFormat = "@"
$ws.Cells.Item(725, 13).Value = 'durante a etapa Check, são obtidas as evidências de auditoria, que consistem em registros, declarações de fato ou outra informação pertinente aos critérios de auditoria, sendo estas verificáveis ou não.'
$ws.Cells.Item(725, 14).NumberFormat = "@"
$ws.Cells.Item(725, 14).Value = 'na etapa Plan, a organização determinará os aspectos ambientais, sendo estes elementos de atividades, produtos ou serviços que interagem ou podem interagir com o meio ambiente, podendo causar impactos ambientais.'
$ws.Cells.Item(725, 15).NumberFormat = "@"
$ws.Cells.Item(725, 15).Value = 'D'
$ws.Cells.Item(725, 16).Value = 0
$ws.Cells.Item(725, 17).Value = 0

# Row 726
$ws.Cells.Item(726, 1).Value = 726
$ws.Cells.Item(726, 2).NumberFormat = "@"
$ws.Cells.Item(726, 2).Value = 'NUCEPE'
$ws.Cells.Item(726, 3).NumberFormat = "@"
$ws.Cells.Item(726, 3).Value = 'PC-PI'
$ws.Cells.Item(726, 4).NumberFormat = "@"
$ws.Cells.Item(726, 4).Value = '2018'
$ws.Cells.Item(726, 5).NumberFormat = "@"
$ws.Cells.Item(726, 5).Value = '“Com a temática da sustentabilidade em alta, nos dias de hoje, fica cada vez mais evidente que a consciência ambiental desempenha um papel definitivo na construção da cidadania. De forma crescente, as pessoas avaliam seus comportamentos em sociedade e como eles se refletem na conservação do nosso ecossistema. Neste sentido o Sistema de Gestão Ambiental (SGA) vem para balizar as ações corporativas em busca do equilíbrio do homem, da indústria e do meio ambiente. O objetivo da Gestão Ambiental é a busca permanente de melhoria da qualidade ambiental dos serviços, produtos e ambiente de trabalho de qualquer organização”., Assinale a alternativa que contém os processos indispensáveis ao SGA.'
$ws.Cells.Item(726, 6).NumberFormat = "@"
$ws.Cells.Item(726, 6).Value = 'Conhecimentos Específicos'
$ws.Cells.Item(726, 7).NumberFormat = "@"
$ws.Cells.Item(726, 7).Value = 'Sustentabilidade'
$ws.Cells.Item(726, 8).NumberFormat = "@"
$ws.Cells.Item(726, 8).Value = 'Médio'
$ws.Cells.Item(726, 9).NumberFormat = "@"
$ws.Cells.Item(726, 9).Value = 'ME'
$ws.Cells.Item(726, 10).NumberFormat = "@"
$ws.Cells.Item(726, 10).Value = 'Realizar internamente uma autoavaliação e demonstração das conformidades.'
$ws.Cells.Item(726, 11).NumberFormat = "@"
$ws.Cells.Item(726, 11).Value = 'Implementar, manter e aprimorar um sistema de gestão ambiental, assegurando-se de sua conformidade com sua política ambiental definida e demonstrando internacionalmente a competência da empresa no quesito produtividade.'
$ws.Cells.Item(726, 12).NumberFormat = "@"
$ws.Cells.Item(726, 12).Value = 'Buscar certificação/registro do seu sistema de gestão ambiental com base na organização exclusivamente interna.'
$ws.Cells.Item(726, 13).NumberFormat = "@"
$ws.Cells.Item(726, 13).Value = 'Implementar um sistema de gestão ambiental; assegurar-se de sua conformidade com a política ambiental do país, definindo metas a partir de referências a processos de gestão de outras organizações ou empresas.'
$ws.Cells.Item(726, 14).NumberFormat = "@"
$ws.Cells.Item(726, 14).Value = 'Implementar, manter e aprimorar um sistema de gestão ambiental e assegurar-se de sua conformidade com a sua política ambiental definida.'
$ws.Cells.Item(726, 15).NumberFormat = "@"
$ws.Cells.Item(726, 15).Value = 'E'
$ws.Cells.Item(726, 16).Value = 0
$ws.Cells.Item(726, 17).Value = 0
$ws.Cells.Item(726, 18).NumberFormat = "@"
$ws.Cells.Item(726, 18).Value = ''
